$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Update the distance-formula documentation: "Manhattan" -> "haversine".
#
# The canonical diff shows the sentence
#   "...destination, using the Manhattan distance formula..."
# becoming
#   "...destination, using the haversine distance formula..."
# and the edit is recorded as three runs (the unchanged prefix, the newly
# typed replacement text, and the unchanged suffix) with the document's
# "_GoBack" bookmark (which marks the site of the most recent edit) moved
# from its old location to right after the newly typed word "haversine".
# We reproduce that exact run layout below.
# ---------------------------------------------------------------------------

$rng = $d.Content
$found = $rng.Find.Execute("destination, using the Manhattan", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target sentence fragment to replace"
}

$editStart = $rng.Start
$newFragment = "destination, using the haversine"

# Retype the located fragment (this mirrors the user selecting the phrase and
# typing the replacement).
$rng.Text = $newFragment

# Force a run boundary right before the freshly typed text by temporarily
# bookmarking that position, then remove the temporary bookmark once the new
# "_GoBack" bookmark has been placed (the split it introduces persists).
$splitRange = $d.Range($editStart, $editStart)
$d.Bookmarks.Add("ZZZ_TEMP_SPLIT", $splitRange)

# Move "_GoBack" so it sits right after the newly typed word "haversine",
# matching where the cursor would be after typing it.
$goBackPos = $editStart + $newFragment.Length
$goBackRange = $d.Range($goBackPos, $goBackPos)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $goBackRange)

# Drop the temporary helper bookmark; the run split it forced remains.
$d.Bookmarks.Item("ZZZ_TEMP_SPLIT").Delete()
